$d = $word.ActiveDocument

# 1. Remove the standalone "Meta description: ..." paragraph that used to
#    follow the title heading (its content is being relocated/rewritten
#    into the final two paragraphs of the document, see step 2).
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
  $para = $d.Paragraphs($i)
  if ($para.Range.Text -like "Meta description*") {
    $para.Range.Delete()
    break
  }
}

# 2. The final paragraph ("Prompt: Create a feature image ...") is replaced
#    by two new paragraphs:
#      - a bold "Play Fruit Combinator Free - Innovative Slot Game" line
#      - an italic meta-description line (without the "Meta description: "
#        label prefix)
$w_ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$newParasXml = '<w:p ' + $w_ns + '><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Fruit Combinator Free - Innovative Slot Game</w:t></w:r></w:p>' + `
  '<w:p ' + $w_ns + '><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Read our review of Fruit Combinator, an innovative online slot game with low volatility. Play it for free on the developer''s website.</w:t></w:r></w:p>'

for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
  $para = $d.Paragraphs($i)
  if ($para.Range.Text -like "Prompt: Create a feature image*") {
    $para.Range.InsertXML($newParasXml) | Out-Null
    break
  }
}
